$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 15
$ws.Range("H15").Value = 1719.8334
$ws.Range("I15").Value = 1719.8334
$ws.Range("K15").Value = 5159.5002
$ws.Range("M15").Value = -4990.5002
# row 86
$ws.Range("H86").Value = 6785.364
$ws.Range("J86").Value = 7273.5
$ws.Range("L86").Value = 7273.5
$ws.Range("N86").Value = -9519.5
# row 89
$ws.Range("H89").Value = 6785.364
$ws.Range("J89").Value = 7273.5
$ws.Range("L89").Value = 36367.5
$ws.Range("N89").Value = -47599.5
# row 116
$ws.Range("H116").Value = 8503.333000000001
$ws.Range("I116").Value = 7130.625
$ws.Range("J116").Value = 11248.75
$ws.Range("K116").Value = 7130.625
$ws.Range("L116").Value = 11248.75
$ws.Range("M116").Value = -3688.625
$ws.Range("N116").Value = -18132.75
# row 131
$ws.Range("H131").Value = 5978.727
$ws.Range("J131").Value = 10544.333
$ws.Range("L131").Value = 31632.999
$ws.Range("N131").Value = -41712.999
# row 137
$ws.Range("H137").Value = 2665.1
$ws.Range("I137").Value = 2527.889
$ws.Range("J137").Value = 3900
$ws.Range("K137").Value = 7583.667
$ws.Range("L137").Value = 11700
$ws.Range("M137").Value = -5033.667
$ws.Range("N137").Value = -16800
# row 138
$ws.Range("H138").Value = 4259.986
$ws.Range("J138").Value = 4953.5186
$ws.Range("L138").Value = 14860.5558
$ws.Range("N138").Value = -25140.5558
# row 141
$ws.Range("H141").Value = 47139.08
$ws.Range("I141").Value = 883.9
$ws.Range("K141").Value = 2651.7
$ws.Range("M141").Value = 2528.3

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 74
$ws.Range("H74").Value = 14613.667
$ws.Range("I74").Value = 3377.75
$ws.Range("J74").Value = 23602.4
$ws.Range("K74").Value = 3377.75
$ws.Range("L74").Value = 23602.4
$ws.Range("M74").Value = -2503.75
$ws.Range("N74").Value = -25350.4
# row 77
$ws.Range("H77").Value = 14613.667
$ws.Range("I77").Value = 3377.75
$ws.Range("J77").Value = 23602.4
$ws.Range("K77").Value = 16888.75
$ws.Range("L77").Value = 118012
$ws.Range("M77").Value = -12520.75
$ws.Range("N77").Value = -126748
# row 97
$ws.Range("H97").Value = 847.2727
$ws.Range("I97").Value = 817.8946999999999
$ws.Range("K97").Value = 817.8946999999999
$ws.Range("M97").Value = -321.8946999999999
# row 112
$ws.Range("H112").Value = 50096.75
$ws.Range("J112").Value = 50096.75
$ws.Range("L112").Value = 50096.75
$ws.Range("N112").Value = -53050.75
# row 122
$ws.Range("H122").Value = 2331.6365
$ws.Range("I122").Value = 2611.75
$ws.Range("J122").Value = 1584.6666
$ws.Range("K122").Value = 7835.25
$ws.Range("L122").Value = 4753.9998
$ws.Range("M122").Value = -5385.25
$ws.Range("N122").Value = -9653.9998

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 62
$ws.Range("H62").Value = 57499.5
$ws.Range("J62").Value = 57499.5
$ws.Range("L62").Value = 57499.5
$ws.Range("N62").Value = -58871.5
# row 65
$ws.Range("H65").Value = 57499.5
$ws.Range("J65").Value = 57499.5
$ws.Range("L65").Value = 172498.5
$ws.Range("N65").Value = -179362.5
# row 75
$ws.Range("H75").Value = 36738
$ws.Range("I75").Value = 36738
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 36738
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -35802
$ws.Range("N75").ClearContents()
# row 78
$ws.Range("H78").Value = 36738
$ws.Range("I78").Value = 36738
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 110214
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = -105534
$ws.Range("N78").Value = -105534
# row 99
$ws.Range("H99").Value = 2939.4827
$ws.Range("I99").Value = 2592.4783
$ws.Range("J99").Value = 4269.6665
$ws.Range("K99").Value = 2592.4783
$ws.Range("L99").Value = 4269.6665
$ws.Range("M99").Value = -1094.4783
$ws.Range("N99").Value = -7265.6665

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 2
$ws.Range("H2").Value = 3500.6667
$ws.Range("I2").Value = 2649.75
$ws.Range("J2").Value = 5202.5
$ws.Range("K2").Value = 2649.75
$ws.Range("L2").Value = 5202.5
$ws.Range("M2").Value = -2536.75
$ws.Range("N2").Value = -5428.5
# row 4
$ws.Range("H4").Value = 9949.5
$ws.Range("I4").Value = 9949.5
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 9949.5
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -9837.5
$ws.Range("N4").ClearContents()
# row 31
$ws.Range("H31").Value = 5616617.5
$ws.Range("I31").Value = 6734339
$ws.Range("J31").Value = 28009.2
$ws.Range("K31").Value = 6734339
$ws.Range("L31").Value = 28009.2
$ws.Range("M31").Value = -6734044
$ws.Range("N31").Value = -28599.2
# row 34
$ws.Range("H34").Value = 5616617.5
$ws.Range("I34").Value = 6734339
$ws.Range("J34").Value = 28009.2
$ws.Range("K34").Value = 6734339
$ws.Range("L34").Value = 28009.2
$ws.Range("M34").Value = -6734137
$ws.Range("N34").Value = -28413.2
# row 94
$ws.Range("H94").Value = 849.619
$ws.Range("J94").Value = 642
$ws.Range("L94").Value = 642
$ws.Range("N94").Value = -1544
# row 115
$ws.Range("H115").Value = 95999
$ws.Range("J115").Value = 95999
$ws.Range("L115").Value = 95999
$ws.Range("N115").Value = -98349
# row 122
$ws.Range("H122").Value = 1753.4
$ws.Range("I122").Value = 1956.1428
$ws.Range("K122").Value = 5868.428400000001
$ws.Range("M122").Value = -3418.428400000001

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 5
$ws.Range("H5").Value = 8663.833000000001
$ws.Range("I5").Value = 396.6
$ws.Range("K5").Value = 1189.8
$ws.Range("M5").Value = -1077.8
# row 46
$ws.Range("H46").Value = 572.125
$ws.Range("I46").Value = 315.4
$ws.Range("K46").Value = 946.1999999999999
$ws.Range("M46").Value = -855.1999999999999
# row 107
$ws.Range("H107").Value = 6536761
$ws.Range("J107").Value = 841.93024
$ws.Range("L107").Value = 2525.79072
$ws.Range("N107").Value = -6365.79072
# row 110
$ws.Range("H110").Value = 11659.667
$ws.Range("I110").Value = 11659.667
$ws.Range("K110").Value = 34979.001
$ws.Range("M110").Value = -30889.001
# row 122
$ws.Range("H122").Value = 2380.25
$ws.Range("I122").Value = 862.2222
$ws.Range("J122").Value = 6934.3335
$ws.Range("K122").Value = 7759.999800000001
$ws.Range("L122").Value = 62409.0015
$ws.Range("M122").Value = -5309.999800000001
$ws.Range("N122").Value = -67309.0015
# row 125
$ws.Range("H125").Value = 17857.143
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 17857.143
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 53571.429
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -63411.429
# row 129
$ws.Range("H129").Value = 24753972
$ws.Range("I129").Value = 33001950
$ws.Range("J129").Value = 10033
$ws.Range("K129").Value = 99005850
$ws.Range("L129").Value = 30099
$ws.Range("M129").Value = -99000850
$ws.Range("N129").Value = -40099
# row 131
$ws.Range("H131").Value = 20653.672
$ws.Range("J131").Value = 3607.8823
$ws.Range("L131").Value = 10823.6469
$ws.Range("N131").Value = -20903.6469
# row 132
$ws.Range("H132").Value = 1724.9375
$ws.Range("I132").Value = 1409.4
$ws.Range("J132").Value = 2250.8333
$ws.Range("K132").Value = 12684.6
$ws.Range("L132").Value = 20257.4997
$ws.Range("M132").Value = -10154.6
$ws.Range("N132").Value = -25317.4997
# row 135
$ws.Range("H135").Value = 8663.833000000001
$ws.Range("I135").Value = 396.6
$ws.Range("K135").Value = 3569.4
$ws.Range("M135").Value = -1034.4

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 102
$ws.Range("H102").Value = 1853.8182
$ws.Range("I102").Value = 1894.4286
$ws.Range("J102").Value = 1001
$ws.Range("K102").Value = 1894.4286
$ws.Range("L102").Value = 1001
$ws.Range("M102").Value = -272.4286
$ws.Range("N102").Value = -4245
# row 126
$ws.Range("H126").Value = 20946.223
$ws.Range("I126").Value = 28594.416
$ws.Range("J126").Value = 5649.8335
$ws.Range("K126").Value = 85783.24800000001
$ws.Range("L126").Value = 16949.5005
$ws.Range("M126").Value = -83313.24800000001
$ws.Range("N126").Value = -21889.5005
# row 132
$ws.Range("H132").Value = 265864.6
$ws.Range("I132").Value = 288370.84
$ws.Range("J132").Value = 3291.3333
$ws.Range("K132").Value = 865112.52
$ws.Range("L132").Value = 9873.999899999999
$ws.Range("M132").Value = -862582.52
$ws.Range("N132").Value = -14933.9999

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 122
$ws.Range("H122").Value = 5981.143
$ws.Range("I122").Value = 2710.4
$ws.Range("K122").Value = 8131.200000000001
$ws.Range("M122").Value = -5681.200000000001
# row 126
$ws.Range("H126").Value = 3779.9546
$ws.Range("I126").Value = 3545.6667
$ws.Range("J126").Value = 4282
$ws.Range("K126").Value = 10637.0001
$ws.Range("L126").Value = 12846
$ws.Range("M126").Value = -8167.000100000001
$ws.Range("N126").Value = -17786
# row 132
$ws.Range("H132").Value = 2706.5264
$ws.Range("I132").Value = 2436.7058
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 7310.117400000001
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -4780.117400000001
